$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete column D ("layer_init_size") - shifts everything right of it one column left
$ws.Range("D1").EntireColumn.Delete() | Out-Null

# After deletion, headers are:
# A Net, B random_state, C num_layers, D layers_size, E net_param, F criterion,
# G learning_rate, H optimizer, I epochs, J vali_best_epoch, K vali_best_acc,
# L vali_best_loss, M vali_best_R, N test_acc, O test_loss, P test_R, Q predicted

# 2. Insert 3 new columns before the current "test_acc" column (N), to hold the
#    new surrogate validation metrics.
$ws.Range("N1:P1").EntireColumn.Insert() | Out-Null

$ws.Range("N1").Value = "vali_sur_acc"
$ws.Range("O1").Value = "vali_sur_loss"
$ws.Range("P1").Value = "vali_sur_R"

# Copy header style from a neighboring header cell so formatting matches.
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1:P1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# 3. Append 2 new columns after "predicted" (now column T) for confusion matrix data.
$ws.Range("U1").Value = "c_matrix"
$ws.Range("V1").Value = "c_matrix_perc"

$ws.Range("T1").Copy() | Out-Null
$ws.Range("U1:V1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the used dimension reference to match the new extent.
$ws.Range("A1:V1").Select() | Out-Null
